$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Title slide text was "Database Design Chat 010"; rename it to
# "Database Design Chat 015". Only the "Chat 010" -> "Chat 015" portion
# of the text actually changed, so replace just that substring; this
# mirrors how the slide was edited (selecting "Chat 010" and retyping
# "Chat 015"), leaving "Database Design " untouched and producing two
# separate runs for the paragraph.
$oldPart = "Chat 010"
$newPart = "Chat 015"
$startPos = $tr.Text.IndexOf($oldPart) + 1
$target = $tr.Characters($startPos, $oldPart.Length)
$target.Text = $newPart
